# Update "想去人数" (interest count) figures on the 展览 and 全部类型 sheets
# to reflect newly generated output, per commit 456a3b4.

$wb = $excel.ActiveWorkbook

# Sheet "展览": row 4 (南宁·2024良牙动漫秋季盛典) and row 6 (南宁·花海演绎二次元水上派对)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F4").Value = 5012
$wsExhibit.Range("F6").Value = 36

# Sheet "全部类型": row 4 (南宁·2024良牙动漫秋季盛典) and row 7 (南宁·花海演绎二次元水上派对)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 5012
$wsAll.Range("F7").Value = 36
